$wb = $excel.ActiveWorkbook

# --- "Estadisticos 2P" (segundo parcial): everyone now accounted for, 0 reprobados ---
$ws2p = $wb.Worksheets.Item("Estadisticos 2P")
$ws2p.Range("D2").Value = 0
$ws2p.Range("E2").Value = 0
$ws2p.Range("F2").Value = 39
$ws2p.Range("G2").Value = 100
$ws2p.Range("H2").Value = 8.800000000000001

# --- "Estadisticos Final": same update, 0 reprobados ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("E2").Value = 0
$wsFinal.Range("F2").Value = 39
$wsFinal.Range("G2").Value = 100
$wsFinal.Range("H2").Value = 9.300000000000001

# --- "Rescatables": the two previously-rescatable students are resolved, remove their rows ---
$wsRescatables = $wb.Worksheets.Item("Rescatables")
$wsRescatables.Rows.Item(2).Delete()
$wsRescatables.Rows.Item(2).Delete()
